$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B5").Value = "мечо е номер едно"
$ws.Range("A5").Value = "ала бала"
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Range("A6").Select() | Out-Null
